$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "stehlen",
    "jubeln",
    "fangen",
    "dienen",
    "quälen",
    "jagen",
    "fallen",
    "machen",
    "bellen",
    "wirken",
    "wüten",
    "feiern",
    "mauern",
    "spinnen",
    "zögern",
    "streichen",
    "grüßen",
    "backen",
    "sperren",
    "sprengen",
    "erben",
    "fließen",
    "fällen",
    "gründen",
    "boxen",
    "dringen",
    "bergen",
    "planen",
    "schmecken",
    "zielen",
    "liefern",
    "seufzen",
    "sorgen",
    "ärgern",
    "graben",
    "decken",
    "bauen",
    "sinken",
    "scheitern",
    "rufen",
    "ehren",
    "töten",
    "pfeifen",
    "stecken",
    "schrecken",
    "formen",
    "kümmern",
    "betteln",
    "filmen",
    "klingen",
    "rasen",
    "arten",
    "scheinen",
    "geben",
    "schwören",
    "werfen",
    "platzen",
    "heilen",
    "biegen",
    "trauen",
    "hauen",
    "siegen",
    "schlucken",
    "treiben",
    "malen",
    "gelten",
    "enden",
    "schulden",
    "bluten",
    "spielen",
    "knarren",
    "fahren",
    "trennen",
    "folgen",
    "klettern",
    "münzen",
    "mögen",
    "kichern",
    "greifen",
    "stammen",
    "warnen",
    "heulen",
    "lügen",
    "spüren",
    "schreiten",
    "lesen",
    "schwingen",
    "helfen",
    "werden",
    "weichen",
    "tollen",
    "altern",
    "wachsen",
    "drehen",
    "flüchten",
    "zeigen",
    "sichern",
    "suchen",
    "äußern",
    "wehtun",
    "loben",
    "schwächen",
    "bitten",
    "kosten",
    "ändern",
    "pflanzen",
    "sterben",
    "heben",
    "schenken",
    "lockern",
    "räumen",
    "brauchen",
    "saufen",
    "achten",
    "freuen",
    "zünden",
    "hören",
    "flehen",
    "tropfen",
    "reizen",
    "wenden",
    "runden",
    "führen",
    "wundern",
    "liegen",
    "irren",
    "kehren",
    "fischen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}

Write-Output "done"
